# Updated cryptos list on Mon Apr  1 02:30:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, $cellRef, $value)
    # Force the cell to hold a text string (not auto-converted to a number)
    # while keeping the cell's original (default) style afterwards.
    $cell = $sheet.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "70.795.33"
$ws.Range("E2").Value = "  +1.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.607.63"
$ws.Range("E3").Value = "  +2.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - Solana
Set-TextValue $ws "D5" "202.96"
$ws.Range("E5").Value = "  +3.22%  "

# Row 6 - BNB
Set-TextValue $ws "D6" "602.95"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.629"
$ws.Range("E7").Value = "  +0.72%  "

# Row 8 - USDC
Set-TextValue $ws "D8" "1.00"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +5.99%  "

# Row 10 - Cardano
Set-TextValue $ws "D10" "0.648"
$ws.Range("E10").Value = "  +0.08%  "

# Row 11 - Avalanche
Set-TextValue $ws "D11" "53.73"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12 - ShibaInu
Set-TextValue $ws "D12" "0.0000303"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +1.57%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.180.85"
$ws.Range("E14").Value = "  +2.25%  "

# Row 15 - BitcoinCash
Set-TextValue $ws "D15" "682.37"
$ws.Range("E15").Value = "  +13.56%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "70.866.13"
$ws.Range("E16").Value = "  +1.23%  "

# Row 17 - Chainlink
Set-TextValue $ws "D17" "19.18"
$ws.Range("E17").Value = "  +0.82%  "

# Row 18 - Uniswap
Set-TextValue $ws "D18" "12.78"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.600.97"
$ws.Range("E19").Value = "  +2.00%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.36%  "

# Row 21 - Polygon
Set-TextValue $ws "D21" "1.00"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextValue $ws "D22" "18.61"
$ws.Range("E22").Value = "  +3.31%  "

# Row 23 - Litecoin
Set-TextValue $ws "D23" "110.90"
$ws.Range("E23").Value = "  +7.11%  "

# Row 24 - Toncoin
Set-TextValue $ws "D24" "5.37"
$ws.Range("E24").Value = "  +3.86%  "

# Row 25 - PancakeSwap
Set-TextValue $ws "D25" "4.64"
$ws.Range("E25").Value = "  -0.04%  "

# Row 26 - ImmutableX
$ws.Range("E26").Value = "  -0.73%  "

# Row 27 - RenderToken
Set-TextValue $ws "D27" "10.59"
$ws.Range("E27").Value = "  -2.05%  "

# Row 28 - LEO
Set-TextValue $ws "D28" "6.01"
$ws.Range("E28").Value = "  -0.71%  "

# Row 29 - Filecoin
Set-TextValue $ws "D29" "10.21"
$ws.Range("E29").Value = "  +6.70%  "

# Row 30 - EthereumClassic
Set-TextValue $ws "D30" "34.47"
$ws.Range("E30").Value = "  +3.38%  "

# Row 31 - dogwifhat
Set-TextValue $ws "D31" "4.52"
$ws.Range("E31").Value = "  +5.86%  "

# Row 32 - NEARProtocol
Set-TextValue $ws "D32" "7.20"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33 - Cosmos
Set-TextValue $ws "D33" "12.25"
$ws.Range("E33").Value = "  -0.87%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -0.14%  "

# Row 35 - OKB
Set-TextValue $ws "D35" "63.66"
$ws.Range("E35").Value = "  -0.16%  "

# Row 36 - Maker
$ws.Range("D36").Value = "3.895.53"
$ws.Range("E36").Value = "  +3.20%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "0.0₃0853"
$ws.Range("E37").Value = "  +4.49%  "

# Row 38 - Dai
Set-TextValue $ws "D38" "0.999"
$ws.Range("E38").Value = "  -0.07%  "

# Row 39 - Bittensor
Set-TextValue $ws "D39" "513.89"
$ws.Range("E39").Value = "  +1.65%  "

# Row 40 - Fetch.AI
Set-TextValue $ws "D40" "3.03"
$ws.Range("E40").Value = "  -4.75%  "

# Row 41 - InjectiveProtocol
Set-TextValue $ws "D41" "37.13"
$ws.Range("E41").Value = "  +1.78%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +0.79%  "

# Row 43 - was Kaspa, now TheGraph
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D43" "0.387"
$ws.Range("E43").Value = "  -1.26%  "

# Row 44 - was TheGraph, now Kaspa
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D44" "0.140"
$ws.Range("E44").Value = "  +4.71%  "

# Row 45 - VeChain
Set-TextValue $ws "D45" "0.0468"
$ws.Range("E45").Value = "  +4.23%  "

# Row 46 - ThetaToken
Set-TextValue $ws "D46" "3.07"
$ws.Range("E46").Value = "  +8.78%  "

# Row 47 - ApeXProtocol
Set-TextValue $ws "D47" "3.42"
$ws.Range("E47").Value = "  +5.49%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +1.05%  "

# Row 49 - THORChain
Set-TextValue $ws "D49" "8.65"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  -0.32%  "

# Row 51 - was CoreDAO, now Jupiter
$ws.Range("B51").Value = "Jupiter"
$ws.Range("C51").Value = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
Set-TextValue $ws "D51" "1.82"
$ws.Range("E51").Value = "  +20.96%  "
